$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2576.2144
$ws.Range("I28").Value = 2551.6365
$ws.Range("J28").Value = 2666.3333
$ws.Range("K28").Value = 2551.6365
$ws.Range("L28").Value = 2666.3333
$ws.Range("M28").Value = -2066.6365
$ws.Range("N28").Value = -3636.3333
$ws.Range("H132").Value = 2002.7693
$ws.Range("I132").Value = 1828.1111
$ws.Range("K132").Value = 5484.3333
$ws.Range("M132").Value = -2954.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2333.9368
$ws.Range("I32").Value = 1431.375
$ws.Range("J32").Value = 11617.429
$ws.Range("K32").Value = 1431.375
$ws.Range("L32").Value = 11617.429
$ws.Range("M32").Value = -1144.375
$ws.Range("N32").Value = -12191.429
$ws.Range("H61").Value = 1413.5946
$ws.Range("I61").Value = 1208.8
$ws.Range("K61").Value = 1208.8
$ws.Range("M61").Value = -996.8
$ws.Range("H122").Value = 3665.8484
$ws.Range("I122").Value = 3118.5386
$ws.Range("K122").Value = 9355.6158
$ws.Range("M122").Value = -6905.6158
$ws.Range("H132").Value = 5794.792
$ws.Range("I132").Value = 3984.855
$ws.Range("K132").Value = 11954.565
$ws.Range("M132").Value = -9424.565000000001
$ws.Range("H136").Value = 1413.5946
$ws.Range("I136").Value = 1208.8
$ws.Range("K136").Value = 3626.4
$ws.Range("M136").Value = -1076.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1346.2727
$ws.Range("I20").Value = 1326.5
$ws.Range("K20").Value = 1326.5
$ws.Range("M20").Value = -1079.5
$ws.Range("H64").Value = 398.25
$ws.Range("J64").Value = 434.42856
$ws.Range("L64").Value = 434.42856
$ws.Range("N64").Value = -884.4285600000001
$ws.Range("H67").Value = 398.25
$ws.Range("J67").Value = 434.42856
$ws.Range("L67").Value = 434.42856
$ws.Range("N67").Value = -1994.42856
$ws.Range("H86").Value = 1986.0769
$ws.Range("J86").Value = 2322.077
$ws.Range("L86").Value = 2322.077
$ws.Range("N86").Value = -4568.077
$ws.Range("H89").Value = 1986.0769
$ws.Range("J89").Value = 2322.077
$ws.Range("L89").Value = 11610.385
$ws.Range("N89").Value = -22842.385
$ws.Range("H99").Value = 58825576
$ws.Range("I99").Value = 90910820
$ws.Range("K99").Value = 90910820
$ws.Range("M99").Value = -90909322
$ws.Range("H105").Value = 4487.2856
$ws.Range("J105").Value = 4999.5
$ws.Range("L105").Value = 4999.5
$ws.Range("N105").Value = -8493.5
$ws.Range("H134").Value = 3153.1538
$ws.Range("I134").Value = 2175.4119
$ws.Range("K134").Value = 6526.2357
$ws.Range("M134").Value = -3991.2357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2111
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 2111
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H132").Value = 1248.174
$ws.Range("I132").Value = 871.8889
$ws.Range("J132").Value = 2602.8
$ws.Range("K132").Value = 2615.6667
$ws.Range("L132").Value = 7808.400000000001
$ws.Range("M132").Value = -85.66670000000022
$ws.Range("N132").Value = -12868.4
$ws.Range("H134").Value = 1290.2174
$ws.Range("I134").Value = 1288.3658
$ws.Range("K134").Value = 3865.0974
$ws.Range("M134").Value = -1330.0974

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 93.25
$ws.Range("I26").Value = 93.25
$ws.Range("K26").Value = 279.75
$ws.Range("M26").Value = 8.25
$ws.Range("H107").Value = 419
$ws.Range("J107").Value = 460.33334
$ws.Range("L107").Value = 1381.00002
$ws.Range("N107").Value = -5221.000019999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2269.0588
$ws.Range("J102").Value = 4466.3335
$ws.Range("L102").Value = 4466.3335
$ws.Range("N102").Value = -7710.3335
$ws.Range("H132").Value = 1769.3334
$ws.Range("I132").Value = 1821.091
$ws.Range("K132").Value = 5463.272999999999
$ws.Range("M132").Value = -2933.272999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 120000
$ws.Range("I33").Value = 120000
$ws.Range("K33").Value = 120000
$ws.Range("M33").Value = -119710
$ws.Range("H46").Value = 1243
$ws.Range("I46").Value = 1243
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1243
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1055
$ws.Range("N46").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H82").Value = 2149.6956
$ws.Range("I82").Value = 1168.5
$ws.Range("J82").Value = 2904.4614
$ws.Range("K82").Value = 1168.5
$ws.Range("L82").Value = 2904.4614
$ws.Range("M82").Value = -807.5
$ws.Range("N82").Value = -3626.4614
$ws.Range("H85").Value = 2149.6956
$ws.Range("I85").Value = 1168.5
$ws.Range("J85").Value = 2904.4614
$ws.Range("K85").Value = 1168.5
$ws.Range("L85").Value = 2904.4614
$ws.Range("M85").Value = 79.5
$ws.Range("N85").Value = -5400.4614
$ws.Range("H100").Value = 2211.875
$ws.Range("I100").Value = 1650
$ws.Range("J100").Value = 2773.75
$ws.Range("K100").Value = 1650
$ws.Range("L100").Value = 2773.75
$ws.Range("M100").Value = -1109
$ws.Range("N100").Value = -3855.75
$ws.Range("H136").Value = 2853.5789
$ws.Range("I136").Value = 2467.6667
$ws.Range("K136").Value = 7403.000100000001
$ws.Range("M136").Value = -4853.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2146.25
$ws.Range("I100").Value = 1417.3334
$ws.Range("K100").Value = 2834.6668
$ws.Range("M100").Value = -2293.6668
$ws.Range("H107").Value = 692.1875
$ws.Range("I107").Value = 506.15384
$ws.Range("K107").Value = 1518.46152
$ws.Range("M107").Value = 401.5384799999999
$ws.Range("H132").Value = 1555.1852
$ws.Range("I132").Value = 1279.24
$ws.Range("K132").Value = 3837.72
$ws.Range("M132").Value = -1307.72
$ws.Range("H136").Value = 5898.147
$ws.Range("I136").Value = 6318.32
$ws.Range("J136").Value = 4731
$ws.Range("K136").Value = 18954.96
$ws.Range("L136").Value = 14193
$ws.Range("M136").Value = -16404.96
$ws.Range("N136").Value = -19293

